# Version 1.1 update:
#  - "This is a main heading" / "Document-1" TOC-style bookmarks get fresh
#    (regenerated) _Toc names.
#  - A new "puzzlesContent" Heading 2 paragraph (with its own _Toc bookmark)
#    is inserted right before the existing "Document-1" Heading 2 paragraph.

$d = $word.ActiveDocument

# --- 1. Rename the bookmark around "This is a main heading  " ---------------
$bmHeading1 = $d.Bookmarks.Item("_Toc16565293719996121707656749")
$h1Start = $bmHeading1.Start
$h1End = $bmHeading1.End
$bmHeading1.Delete()
$h1Range = $d.Range($h1Start, $h1End)
$d.Bookmarks.Add("_Toc16565903375248968875675980", $h1Range) | Out-Null

# --- 2. Find the "Document-1" heading paragraph and insert a new paragraph
#        ("puzzlesContent  ") right before it, reusing the Heading 2 style. --
$doc1Bookmark = $d.Bookmarks.Item("_Toc16565293720296803199724611")
$doc1ParaRange = $doc1Bookmark.Range
$doc1Para = $doc1ParaRange.Paragraphs.Item(1)
$doc1Index = $doc1Para.Index

$doc1Para.Range.InsertParagraphBefore()

# The freshly inserted (still empty) paragraph now sits at the old index;
# "Document-1" shifted one paragraph further down.
$newPara = $d.Paragraphs.Item($doc1Index)
$newRange = $newPara.Range
$newRange.InsertBefore("puzzlesContent  ")

# --- 3. Add the bookmark for the new "puzzlesContent" paragraph -------------
$puzzlesRange = $d.Range($newRange.Start, $newRange.Start + "puzzlesContent  ".Length)
$d.Bookmarks.Add("_Toc16565903375517092282959150", $puzzlesRange) | Out-Null

# --- 4. Rename the bookmark that now sits on "Document-1" -------------------
$bmDoc1 = $d.Bookmarks.Item("_Toc16565293720296803199724611")
$doc1Start = $bmDoc1.Start
$doc1End = $bmDoc1.End
$bmDoc1.Delete()
$doc1Range = $d.Range($doc1Start, $doc1End)
$d.Bookmarks.Add("_Toc16565903375796663920027541", $doc1Range) | Out-Null
